$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new column before column A: shifts the existing Supplier /
# Contact Name / Email headers (and the trailing blank cell) one column
# to the right, preserving their shared-string values, styles, and
# column widths automatically.
$ws.Columns.Item(1).Insert() | Out-Null

# The new first column holds the "Supplier Code" header.
$ws.Range("A1").Value = "Supplier Code"

# Match the bold/filled header style used by the other header cells.
$ws.Range("B1").Copy() | Out-Null
$ws.Range("A1").PasteSpecial(-4122) | Out-Null

# Give the new column its own width (distinct from the old column-A
# width, which now belongs to column B).
$ws.Columns.Item(1).ColumnWidth = 19

# Move the active selection, matching the saved workbook.
$ws.Range("B4").Select() | Out-Null
